$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 911
$ws.Range("I33").Value = 978.2308
$ws.Range("K33").Value = 978.2308
$ws.Range("M33").Value = -749.2308

$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()

$ws.Range("H53").Value = 235.11765
$ws.Range("I53").Value = 258
$ws.Range("K53").Value = 258
$ws.Range("M53").Value = 379

$ws.Range("H121").Value = 5314.143
$ws.Range("J121").Value = 5314.143
$ws.Range("L121").Value = 15942.429
$ws.Range("N121").Value = -19436.429

$ws.Range("H137").Value = 2240.6938
$ws.Range("I137").Value = 1892.6666
$ws.Range("J137").Value = 2790.2104
$ws.Range("K137").Value = 5677.9998
$ws.Range("L137").Value = 8370.6312
$ws.Range("M137").Value = -3127.9998
$ws.Range("N137").Value = -13470.6312

$ws.Range("H138").Value = 1256193.6
$ws.Range("J138").Value = 1455860.9
$ws.Range("L138").Value = 4367582.699999999
$ws.Range("N138").Value = -4377862.699999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3989.855
$ws.Range("I32").Value = 3295.9395
$ws.Range("K32").Value = 3295.9395
$ws.Range("M32").Value = -3008.9395

$ws.Range("H132").Value = 2399359
$ws.Range("I132").Value = 6581069
$ws.Range("J132").Value = 9810.5
$ws.Range("K132").Value = 19743207
$ws.Range("L132").Value = 29431.5
$ws.Range("M132").Value = -19740677
$ws.Range("N132").Value = -34491.5

$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 7498.3335
$ws.Range("I35").Value = 7498.3335
$ws.Range("K35").Value = 7498.3335
$ws.Range("M35").Value = -7188.3335

$ws.Range("H105").Value = 2226334
$ws.Range("I105").Value = 5001472
$ws.Range("J105").Value = 6223.4
$ws.Range("K105").Value = 5001472
$ws.Range("L105").Value = 6223.4
$ws.Range("M105").Value = -4999725
$ws.Range("N105").Value = -9717.4

$ws.Range("H134").Value = 9970.333000000001
$ws.Range("I134").Value = 4345
$ws.Range("J134").Value = 10673.5
$ws.Range("K134").Value = 13035
$ws.Range("L134").Value = 32020.5
$ws.Range("M134").Value = -10500
$ws.Range("N134").Value = -37090.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 262.75
$ws.Range("I22").Value = 262.75
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 262.75
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 87.25
$ws.Range("N22").ClearContents()

$ws.Range("H31").Value = 5698.644
$ws.Range("J31").Value = 7127.4634
$ws.Range("L31").Value = 7127.4634
$ws.Range("N31").Value = -7717.4634

$ws.Range("H34").Value = 5698.644
$ws.Range("J34").Value = 7127.4634
$ws.Range("L34").Value = 7127.4634
$ws.Range("N34").Value = -7531.4634

$ws.Range("H105").Value = 3573141.5
$ws.Range("I105").Value = 3969046
$ws.Range("J105").Value = 9999
$ws.Range("K105").Value = 3969046
$ws.Range("L105").Value = 9999
$ws.Range("M105").Value = -3967299
$ws.Range("N105").Value = -13493

$ws.Range("H135").Value = 94500
$ws.Range("J135").Value = 94500
$ws.Range("L135").Value = 94500
$ws.Range("N135").Value = -104640

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 412.33334
$ws.Range("I18").Value = 471.3
$ws.Range("J18").Value = 117.5
$ws.Range("K18").Value = 1413.9
$ws.Range("L18").Value = 352.5
$ws.Range("M18").Value = -1244.9
$ws.Range("N18").Value = -690.5

$ws.Range("H68").Value = 2904.5881
$ws.Range("I68").Value = 2412.25
$ws.Range("J68").Value = 2996.186
$ws.Range("K68").Value = 7236.75
$ws.Range("L68").Value = 8988.558000000001
$ws.Range("M68").Value = -6425.75
$ws.Range("N68").Value = -10610.558

$ws.Range("H71").Value = 2904.5881
$ws.Range("I71").Value = 2412.25
$ws.Range("J71").Value = 2996.186
$ws.Range("K71").Value = 21710.25
$ws.Range("L71").Value = 26965.674
$ws.Range("M71").Value = -17654.25
$ws.Range("N71").Value = -35077.674

$ws.Range("H74").Value = 2996.75
$ws.Range("J74").Value = 2996.75
$ws.Range("L74").Value = 8990.25
$ws.Range("N74").Value = -11112.25

$ws.Range("H77").Value = 2996.75
$ws.Range("J77").Value = 2996.75
$ws.Range("L77").Value = 26970.75
$ws.Range("N77").Value = -37578.75

$ws.Range("H87").Value = 125004630
$ws.Range("I87").Value = 166668830
$ws.Range("J87").Value = 11998
$ws.Range("K87").Value = 500006490
$ws.Range("L87").Value = 35994
$ws.Range("M87").Value = -500005242
$ws.Range("N87").Value = -38490

$ws.Range("H90").Value = 125004630
$ws.Range("I90").Value = 166668830
$ws.Range("J90").Value = 11998
$ws.Range("K90").Value = 1500019470
$ws.Range("L90").Value = 107982
$ws.Range("M90").Value = -1500013230
$ws.Range("N90").Value = -120462

$ws.Range("H107").Value = 22222572
$ws.Range("J107").Value = 100000500
$ws.Range("L107").Value = 300001500
$ws.Range("N107").Value = -300005340

$ws.Range("H131").Value = 1686.3077
$ws.Range("I131").Value = 910
$ws.Range("K131").Value = 2730
$ws.Range("M131").Value = 2310

$ws.Range("H132").Value = 12683.9375
$ws.Range("J132").Value = 15599.3
$ws.Range("L132").Value = 140393.7
$ws.Range("N132").Value = -145453.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 11332.667
$ws.Range("J70").Value = 11714.143
$ws.Range("L70").Value = 11714.143
$ws.Range("N70").Value = -12254.143

$ws.Range("H73").Value = 11332.667
$ws.Range("J73").Value = 11714.143
$ws.Range("L73").Value = 11714.143
$ws.Range("N73").Value = -13586.143

$ws.Range("H80").Value = 4392.8
$ws.Range("I80").Value = 4371.625
$ws.Range("K80").Value = 4371.625
$ws.Range("M80").Value = -3373.625

$ws.Range("H83").Value = 4392.8
$ws.Range("I83").Value = 4371.625
$ws.Range("K83").Value = 21858.125
$ws.Range("M83").Value = -16866.125

$ws.Range("H93").Value = 52500
$ws.Range("J93").Value = 15000
$ws.Range("L93").Value = 15000
$ws.Range("N93").Value = -18744

$ws.Range("H97").Value = 1605.5
$ws.Range("I97").Value = 1655.55
$ws.Range("K97").Value = 1655.55
$ws.Range("M97").Value = -1159.55

$ws.Range("H122").Value = 83321.16
$ws.Range("I122").Value = 204440.8
$ws.Range("K122").Value = 613322.3999999999
$ws.Range("M122").Value = -610872.3999999999

$ws.Range("H126").Value = 7581.9443
$ws.Range("I126").Value = 6952.4546
$ws.Range("J126").Value = 8571.143
$ws.Range("K126").Value = 20857.3638
$ws.Range("L126").Value = 25713.429
$ws.Range("M126").Value = -18387.3638
$ws.Range("N126").Value = -30653.429

$ws.Range("H132").Value = 3767.5898
$ws.Range("I132").Value = 2359.5
$ws.Range("J132").Value = 6583.769
$ws.Range("K132").Value = 7078.5
$ws.Range("L132").Value = 19751.307
$ws.Range("M132").Value = -4548.5
$ws.Range("N132").Value = -24811.307

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2979.8
$ws.Range("I68").Value = 2974.75
$ws.Range("K68").Value = 2974.75
$ws.Range("M68").Value = -2225.75

$ws.Range("H71").Value = 2979.8
$ws.Range("I71").Value = 2974.75
$ws.Range("K71").Value = 14873.75
$ws.Range("M71").Value = -11129.75

$ws.Range("H122").Value = 6783.222
$ws.Range("J122").Value = 7244.1665
$ws.Range("L122").Value = 21732.4995
$ws.Range("N122").Value = -26632.4995

$ws.Range("H132").Value = 16136432
$ws.Range("I132").Value = 31253246
$ws.Range("K132").Value = 93759738
$ws.Range("M132").Value = -93757208

$ws.Range("H136").Value = 13572.389
$ws.Range("I136").Value = 9200.6
$ws.Range("J136").Value = 15253.846
$ws.Range("K136").Value = 27601.8
$ws.Range("L136").Value = 45761.538
$ws.Range("M136").Value = -25051.8
$ws.Range("N136").Value = -50861.538

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 19099886
$ws.Range("J81").Value = 28584686
$ws.Range("L81").Value = 57169372
$ws.Range("N81").Value = -57171494

$ws.Range("H84").Value = 19099886
$ws.Range("J84").Value = 28584686
$ws.Range("L84").Value = 285846860
$ws.Range("N84").Value = -285857468

$ws.Range("H122").Value = 3298.8
$ws.Range("J122").Value = 5319
$ws.Range("L122").Value = 15957
$ws.Range("N122").Value = -20857

$ws.Range("H132").Value = 17246.037
$ws.Range("I132").Value = 9036.022000000001
$ws.Range("J132").Value = 58296.11
$ws.Range("K132").Value = 27108.066
$ws.Range("L132").Value = 174888.33
$ws.Range("M132").Value = -24578.066
$ws.Range("N132").Value = -179948.33

$ws.Range("H136").Value = 58828544
$ws.Range("I136").Value = 76925560
$ws.Range("K136").Value = 230776680
$ws.Range("M136").Value = -230776680
